$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string "Total Middle East" block (rows 3686-3740) follows the same
# per-country layout as the rest of the sheet (55 yearly rows, 1965-2019).
# Seed formatting for the new rows by copying the format of the last existing
# data row (row 3684, which already carries the plain (non-bold) number style).
$ws.Range("A3684:C3684").Copy($ws.Range("A3686:C3740"))

# The previous last row (3685) loses its "end of block" bold styling now that
# it is followed by more data -- match it to the plain style used elsewhere.
$ws.Range("C3684").Copy()
$ws.Range("C3685").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 3686; Date = 23743; Country = "Total Middle East"; Barrels = 870.00565430442828 },
    @{ Row = 3687; Date = 24108; Country = "Total Middle East"; Barrels = 894.34242038045772 },
    @{ Row = 3688; Date = 24473; Country = "Total Middle East"; Barrels = 922.96024422263201 },
    @{ Row = 3689; Date = 24838; Country = "Total Middle East"; Barrels = 951.04871811327439 },
    @{ Row = 3690; Date = 25204; Country = "Total Middle East"; Barrels = 984.84839092177651 },
    @{ Row = 3691; Date = 25569; Country = "Total Middle East"; Barrels = 1043.988224157825 },
    @{ Row = 3692; Date = 25934; Country = "Total Middle East"; Barrels = 1098.2894873142525 },
    @{ Row = 3693; Date = 26299; Country = "Total Middle East"; Barrels = 1180.1757072236194 },
    @{ Row = 3694; Date = 26665; Country = "Total Middle East"; Barrels = 1275.9133848665967 },
    @{ Row = 3695; Date = 27030; Country = "Total Middle East"; Barrels = 1360.6817723152933 },
    @{ Row = 3696; Date = 27395; Country = "Total Middle East"; Barrels = 1316.562310673876 },
    @{ Row = 3697; Date = 27760; Country = "Total Middle East"; Barrels = 1501.6015211992765 },
    @{ Row = 3698; Date = 28126; Country = "Total Middle East"; Barrels = 1706.1689259447928 },
    @{ Row = 3699; Date = 28491; Country = "Total Middle East"; Barrels = 1775.7926117267871 },
    @{ Row = 3700; Date = 28856; Country = "Total Middle East"; Barrels = 2018.3628073514913 },
    @{ Row = 3701; Date = 29221; Country = "Total Middle East"; Barrels = 1932.9859880279 },
    @{ Row = 3702; Date = 29587; Country = "Total Middle East"; Barrels = 2116.0385322137345 },
    @{ Row = 3703; Date = 29952; Country = "Total Middle East"; Barrels = 2326.6001881268066 },
    @{ Row = 3704; Date = 30317; Country = "Total Middle East"; Barrels = 2596.5719630330454 },
    @{ Row = 3705; Date = 30682; Country = "Total Middle East"; Barrels = 2832.5710500737168 },
    @{ Row = 3706; Date = 31048; Country = "Total Middle East"; Barrels = 3011.699452196151 },
    @{ Row = 3707; Date = 31413; Country = "Total Middle East"; Barrels = 3001.094099927795 },
    @{ Row = 3708; Date = 31778; Country = "Total Middle East"; Barrels = 3177.4474363134436 },
    @{ Row = 3709; Date = 32143; Country = "Total Middle East"; Barrels = 3309.8131328694667 },
    @{ Row = 3710; Date = 32509; Country = "Total Middle East"; Barrels = 3420.1642871238005 },
    @{ Row = 3711; Date = 32874; Country = "Total Middle East"; Barrels = 3488.553870355991 },
    @{ Row = 3712; Date = 33239; Country = "Total Middle East"; Barrels = 3645.58084058739 },
    @{ Row = 3713; Date = 33604; Country = "Total Middle East"; Barrels = 3851.3607971844503 },
    @{ Row = 3714; Date = 33970; Country = "Total Middle East"; Barrels = 4144.7359845345363 },
    @{ Row = 3715; Date = 34335; Country = "Total Middle East"; Barrels = 4584.5956762283886 },
    @{ Row = 3716; Date = 34700; Country = "Total Middle East"; Barrels = 4600.6411443177631 },
    @{ Row = 3717; Date = 35065; Country = "Total Middle East"; Barrels = 4718.6095407627363 },
    @{ Row = 3718; Date = 35431; Country = "Total Middle East"; Barrels = 4951.4200096755021 },
    @{ Row = 3719; Date = 35796; Country = "Total Middle East"; Barrels = 4861.9457530819218 },
    @{ Row = 3720; Date = 36161; Country = "Total Middle East"; Barrels = 4852.9825312957773 },
    @{ Row = 3721; Date = 36526; Country = "Total Middle East"; Barrels = 5087.4057996598694 },
    @{ Row = 3722; Date = 36892; Country = "Total Middle East"; Barrels = 5320.2230878241608 },
    @{ Row = 3723; Date = 37257; Country = "Total Middle East"; Barrels = 5455.5655945716799 },
    @{ Row = 3724; Date = 37622; Country = "Total Middle East"; Barrels = 5673.7081739788628 },
    @{ Row = 3725; Date = 37987; Country = "Total Middle East"; Barrels = 5997.9598858199779 },
    @{ Row = 3726; Date = 38353; Country = "Total Middle East"; Barrels = 6452.0073769819637 },
    @{ Row = 3727; Date = 38718; Country = "Total Middle East"; Barrels = 6721.9990333413016 },
    @{ Row = 3728; Date = 39083; Country = "Total Middle East"; Barrels = 6972.7757976428093 },
    @{ Row = 3729; Date = 39448; Country = "Total Middle East"; Barrels = 7391.2971528604439 },
    @{ Row = 3730; Date = 39814; Country = "Total Middle East"; Barrels = 7730.3820363724099 },
    @{ Row = 3731; Date = 40179; Country = "Total Middle East"; Barrels = 7987.0544942342603 },
    @{ Row = 3732; Date = 40544; Country = "Total Middle East"; Barrels = 8292.6553892270822 },
    @{ Row = 3733; Date = 40909; Country = "Total Middle East"; Barrels = 8641.9797228655261 },
    @{ Row = 3734; Date = 41275; Country = "Total Middle East"; Barrels = 8867.9129499334595 },
    @{ Row = 3735; Date = 41640; Country = "Total Middle East"; Barrels = 8993.0580440155791 },
    @{ Row = 3736; Date = 42005; Country = "Total Middle East"; Barrels = 8987.209124027524 },
    @{ Row = 3737; Date = 42370; Country = "Total Middle East"; Barrels = 9190.6470876332005 },
    @{ Row = 3738; Date = 42736; Country = "Total Middle East"; Barrels = 9156.4286395324998 },
    @{ Row = 3739; Date = 43101; Country = "Total Middle East"; Barrels = 9173.7268183490305 },
    @{ Row = 3740; Date = 43466; Country = "Total Middle East"; Barrels = 9415.7213340026901 }
)

foreach ($row in $newRows) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Date
    $ws.Cells.Item($row.Row, 2).Value = $row.Country
    $ws.Cells.Item($row.Row, 3).Value = $row.Barrels
}

# Match the author's final on-screen selection after auto-saving.
$ws.Range("G3606").Select()
